$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 3.9
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("AB3").Value = 8.5
$ws.Range("AL3").Value = 19
$ws.Range("AO3").Value = 41
$ws.Range("AP3").Value = 51

# Row 5 updates
$ws.Range("G5").Value = 2.05
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 1.5
$ws.Range("AR5").Value = 4.2
$ws.Range("AS5").Value = 1.22

# Row 6 updates
$ws.Range("G6").Value = 1.69
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.4
$ws.Range("K6").Value = 2.2
$ws.Range("Y6").Value = 1.83
$ws.Range("Z6").Value = 1.83
$ws.Range("AB6").Value = 8
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 13
$ws.Range("AH6").Value = 6.5
$ws.Range("AM6").Value = 17

# Row 8 updates
$ws.Range("I8").Value = 3.6
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 1.53
$ws.Range("P8").Value = 2.38
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("S8").Value = 2.6
$ws.Range("T8").Value = 1.48
$ws.Range("AB8").Value = 8
$ws.Range("AG8").Value = 6.5
$ws.Range("AH8").Value = 6.5
